# IPA and ENW login page and Xls fixes
# Adds a new test-case row (row 12) to the "Test Cases" sheet describing
# the STeAM <-> Facebook account merge / TRUID verification scenario,
# and updates the sheet view / selection to point at the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New row of data (row 12), mirroring the style/format of the existing rows.
$ws.Range("A12").Value = "IPA054"
$ws.Range("B12").Value = "OPQA-4233 || OPQA-4271"
$ws.Range("C12").Value = " Verify that the system is able to merge New STeAM account and Activated Facebook account and after merge verify STeAM TRUID is changed"
$ws.Range("D12").Value = "Y"
$ws.Range("E12").Value = ""

# Match formatting used by the sibling rows (border + fill/wrap as appropriate).
$ws.Range("A12").Style = $ws.Range("A11").Style
$ws.Range("B12").Style = $ws.Range("A11").Style
$ws.Range("C12").Style = $ws.Range("C9").Style
$ws.Range("D12").Style = $ws.Range("D9").Style
$ws.Range("E12").Style = $ws.Range("A11").Style

$ws.Rows.Item(12).RowHeight = 45

# Update the view so the new row is visible and selected, like after
# entering data into the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("A12:E12").Select()
